$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 9943.8799999999992
$ws.Range("B4").Value = 9927
$ws.Range("C4").Value = 309.02999999999997
$ws.Range("D4").Value = 309.55
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.17
$ws.Range("G4").Value = 42608.637835648151
$ws.Range("H4").Value = $true

# Copy the date/time number formatting from G3 (row above) onto G4 so it
# keeps using the same style index (style 1 -> numFmtId 22 date format).
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
